$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new weekly price record for Perejil (Vega Modelo de Temuco) was added on
# 2022-07-04 (serial 44746). It belongs right after the existing record at
# row 304, so insert a fresh row there and push the rest of the table
# (old rows 304-324) down by one - old row 324 ends up at row 325.
$ws.Rows.Item(304).Insert()

$ws.Cells.Item(304, 1).Value = 10
$ws.Cells.Item(304, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(304, 3).Value = "La Araucanía"
$ws.Cells.Item(304, 4).Value = 44746
$ws.Cells.Item(304, 5).Value = 9
$ws.Cells.Item(304, 6).Value = 100112044
$ws.Cells.Item(304, 7).Value = "Perejil"
$ws.Cells.Item(304, 8).Value = "Sin especificar"
$ws.Cells.Item(304, 9).Value = "Primera"
$ws.Cells.Item(304, 10).Value = 55
$ws.Cells.Item(304, 11).Value = 4000
$ws.Cells.Item(304, 12).Value = 4000
$ws.Cells.Item(304, 13).Value = 4000
$ws.Cells.Item(304, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(304, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(304, 16).Value = 1333
$ws.Cells.Item(304, 17).Value = 3
$ws.Cells.Item(304, 18).Value = "Hortaliza"
